# Cryptos list refresh - Tue Jan 23 11:42:11 UTC 2024 (GitHub Actions bot).
#
# Columns D (Price) and E (Volume(1h)) are stored as plain TEXT in this sheet
# (e.g. "27.84", "  -4.81%  "), not as numbers. A straight `.Value = "27.82"`
# assignment would let Excel auto-coerce a clean numeric-looking string into a
# real number, which would not match the source data. To keep such D-column
# values as text (without leaving a stray NumberFormat style behind), each of
# those cells is round-tripped through a text NumberFormat and then restored
# to the "Normal" style, which keeps the default (unstyled) cell format.
# E-column values keep their leading/trailing spaces and "%" sign so they are
# never number-like and can be assigned directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '38.895.55'
$ws.Range("E2").Value = '  -4.77%  '
$ws.Range("D3").Value = '2.223.07'
$ws.Range("E3").Value = '  -6.91%  '
$ws.Range("E4").Value = '  +0.04%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '296.46'
$c.Style = "Normal"

$ws.Range("E5").Value = '  -5.48%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '79.65'
$c.Style = "Normal"

$ws.Range("E6").Value = '  -9.83%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.506'
$c.Style = "Normal"

$ws.Range("E7").Value = '  -4.74%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("E9").Value = '  -7.54%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.0771'
$c.Style = "Normal"

$ws.Range("E10").Value = '  -6.43%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '27.82'
$c.Style = "Normal"

$ws.Range("E11").Value = '  -10.97%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '46.15'
$c.Style = "Normal"

$ws.Range("E12").Value = '  -13.12%  '
$ws.Range("E13").Value = '  -1.80%  '
$ws.Range("D14").Value = '2.567.78'
$ws.Range("E14").Value = '  -6.90%  '
$ws.Range("E15").Value = '  -7.48%  '
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '13.95'
$c.Style = "Normal"

$ws.Range("E16").Value = '  -8.02%  '
$ws.Range("D17").Value = '2.232.48'
$ws.Range("E17").Value = '  -5.60%  '
$ws.Range("E18").Value = '  -6.81%  '
$ws.Range("D19").Value = '38.838.34'
$ws.Range("E19").Value = '  -4.79%  '
$ws.Range("D20").Value = '0.0₃0858'
$ws.Range("E20").Value = '  -5.95%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '5.72'
$c.Style = "Normal"

$ws.Range("E21").Value = '  -7.83%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '64.68'
$c.Style = "Normal"

$ws.Range("E22").Value = '  -6.54%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '9.82'
$c.Style = "Normal"

$ws.Range("E23").Value = '  -9.41%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '224.69'
$c.Style = "Normal"

$ws.Range("E24").Value = '  -3.70%  '
$ws.Range("E25").Value = '  -0.15%  '
$ws.Range("E26").Value = '  -10.15%  '
$ws.Range("E27").Value = '  -6.25%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '2.18'
$c.Style = "Normal"

$ws.Range("E28").Value = '  -1.10%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '22.12'
$c.Style = "Normal"

$ws.Range("E29").Value = '  -7.16%  '
$ws.Range("E30").Value = '  -5.76%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '148.62'
$c.Style = "Normal"

$ws.Range("E31").Value = '  -4.83%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '31.01'
$c.Style = "Normal"

$ws.Range("E32").Value = '  -8.85%  '
$ws.Range("E33").Value = '  -0.11%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '4.77'
$c.Style = "Normal"

$ws.Range("E34").Value = '  -8.61%  '
$ws.Range("E35").Value = '  -4.21%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '0.0682'
$c.Style = "Normal"

$ws.Range("E36").Value = '  -6.92%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '0.109'
$c.Style = "Normal"

$ws.Range("E37").Value = '  -4.41%  '
$ws.Range("E38").Value = '  -5.84%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.0952'
$c.Style = "Normal"

$ws.Range("E39").Value = '  -5.25%  '
$ws.Range("B40").Value = 'Celestia'
$ws.Range("C40").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '14.33'
$c.Style = "Normal"

$ws.Range("E40").Value = '  -11.63%  '
$ws.Range("B41").Value = 'ARBITRUM'
$ws.Range("C41").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '1.58'
$c.Style = "Normal"

$ws.Range("E41").Value = '  -8.83%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '3.58'
$c.Style = "Normal"

$ws.Range("E42").Value = '  -6.16%  '
$ws.Range("D43").Value = '1.898.65'
$ws.Range("E43").Value = '  -3.01%  '
$ws.Range("E44").Value = '  -9.66%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.0253'
$c.Style = "Normal"

$ws.Range("E45").Value = '  -6.92%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '16.15'
$c.Style = "Normal"

$ws.Range("E46").Value = '  -7.95%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '8.95'
$c.Style = "Normal"

$ws.Range("E47").Value = '  -4.67%  '
$ws.Range("E48").Value = '  -10.95%  '
$ws.Range("D49").Value = '2.434.43'
$ws.Range("E49").Value = '  -6.94%  '
$ws.Range("E50").Value = '  -6.50%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '87.17'
$c.Style = "Normal"

$ws.Range("E51").Value = '  -7.17%  '
